$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.252.89'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.70%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.618.87'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.11%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '569.15'
$ws.Range('D5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.05%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('E8').Value = '  +4.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.637.28'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.77'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('E11').Value = '  +3.50%  '
$ws.Range('E12').Value = '  +9.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.344'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.080.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '60.184.96'
$ws.Range('D15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '22.10'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +6.37%  '
$ws.Range('E17').Value = '  +2.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.640.37'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.93%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '340.71'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.10%  '
$ws.Range('E21').Value = '  +2.83%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.35'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.75%  '
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.449'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.90%  '
$ws.Range('E26').Value = '  +3.33%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.17%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.71%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0794'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +7.64%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.11%  '
$ws.Range('E31').Value = '  +3.60%  '
$ws.Range('E32').Value = '  +3.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '159.73'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '19.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.49%  '
$ws.Range('E36').Value = '  +4.00%  '
$ws.Range('B37').Value = 'Fetch.AI'
$ws.Range('C37').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.888'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.76%  '
$ws.Range('B38').Value = 'SuiNetwork'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.881'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.82%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('E40').Value = '  +6.13%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '297.03'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.98%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.65'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.11%  '
$ws.Range('E43').Value = '  -0.21%  '
$ws.Range('E44').Value = '  +4.07%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.601'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.37%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0541'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.24%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.78%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '10.68'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '127.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +15.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0235'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.960.67'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.10%  '
